{"js": "// Replace the 100 arithmetic-problem answers in the single 20x5 table,\n// in row-major (top-to-bottom, left-to-right) order, matching the\n// author's commit. Each cell keeps its existing run formatting because\n// we only change the table's cell text values, not the run properties.\n\nconst newValues = [\n  [\"11+70=81\", \"47-8=39\", \"71-24=47\", \"84-57=27\", \"87-76=11\"],\n  [\"45-2=43\", \"69-46=23\", \"61-13=48\", \"65-31=34\", \"49+6=55\"],\n  [\"72-62=10\", \"40+42=82\", \"58+3=61\", \"87+8=95\", \"42+16=58\"],\n  [\"13+20=33\", \"48+2=50\", \"16-4=12\", \"26+28=54\", \"86-72=14\"],\n  [\"51-43=8\", \"76+13=89\", \"44-6=38\", \"97-37=60\", \"40-7=33\"],\n  [\"65-26=39\", \"89-14=75\", \"54-6=48\", \"22+26=48\", \"84-56=28\"],\n  [\"10+74=84\", \"89-56=33\", \"43-31=12\", \"31+58=89\", \"23+59=82\"],\n  [\"43+54=97\", \"44+12=56\", \"62-55=7\", \"18+22=40\", \"60+11=71\"],\n  [\"79-56=23\", \"99-27=72\", \"94-11=83\", \"78+9=87\", \"73-54=19\"],\n  [\"3+8=11\", \"87-53=34\", \"3+28=31\", \"35-2=33\", \"98-11=87\"],\n  [\"57+3=60\", \"59-53=6\", \"99-68=31\", \"58-5=53\", \"80-32=48\"],\n  [\"91-87=4\", \"4+26=30\", \"3+11=14\", \"25+12=37\", \"64+12=76\"],\n  [\"49-28=21\", \"58-10=48\", \"85-2=83\", \"20+5=25\", \"30+6=36\"],\n  [\"63+30=93\", \"55+34=89\", \"23+30=53\", \"75-17=58\", \"47-30=17\"],\n  [\"33+28=61\", \"72+13=85\", \"31+40=71\", \"7+58=65\", \"85+2=87\"],\n  [\"77-18=59\", \"37-9=28\", \"99-92=7\", \"84-25=59\", \"77+10=87\"],\n  [\"64+18=82\", \"50+19=69\", \"38+2=40\", \"96-40=56\", \"79-59=20\"],\n  [\"81-30=51\", \"81-48=33\", \"84-2=82\", \"27-16=11\", \"43+1=44\"],\n  [\"48+34=82\", \"8+49=57\", \"35+29=64\", \"17+49=66\", \"69-55=14\"],\n  [\"23+61=84\", \"77+14=91\", \"38+14=52\", \"95-49=46\", \"34+55=89\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst current = table.values;\nfor (let r = 0; r < current.length && r < newValues.length; r++) {\n  for (let c = 0; c < current[r].length && c < newValues[r].length; c++) {\n    current[r][c] = newValues[r][c];\n  }\n}\ntable.values = current;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem answers in the single 20x5 table,\n# in row-major (top-to-bottom, left-to-right) order, matching the\n# author's commit. Setting Cell(...).Range.Text keeps the cell's\n# end-of-cell marker and the run's existing formatting (font/size),\n# only the visible text content changes.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"11+70=81\", \"47-8=39\", \"71-24=47\", \"84-57=27\", \"87-76=11\"),\n  @(\"45-2=43\", \"69-46=23\", \"61-13=48\", \"65-31=34\", \"49+6=55\"),\n  @(\"72-62=10\", \"40+42=82\", \"58+3=61\", \"87+8=95\", \"42+16=58\"),\n  @(\"13+20=33\", \"48+2=50\", \"16-4=12\", \"26+28=54\", \"86-72=14\"),\n  @(\"51-43=8\", \"76+13=89\", \"44-6=38\", \"97-37=60\", \"40-7=33\"),\n  @(\"65-26=39\", \"89-14=75\", \"54-6=48\", \"22+26=48\", \"84-56=28\"),\n  @(\"10+74=84\", \"89-56=33\", \"43-31=12\", \"31+58=89\", \"23+59=82\"),\n  @(\"43+54=97\", \"44+12=56\", \"62-55=7\", \"18+22=40\", \"60+11=71\"),\n  @(\"79-56=23\", \"99-27=72\", \"94-11=83\", \"78+9=87\", \"73-54=19\"),\n  @(\"3+8=11\", \"87-53=34\", \"3+28=31\", \"35-2=33\", \"98-11=87\"),\n  @(\"57+3=60\", \"59-53=6\", \"99-68=31\", \"58-5=53\", \"80-32=48\"),\n  @(\"91-87=4\", \"4+26=30\", \"3+11=14\", \"25+12=37\", \"64+12=76\"),\n  @(\"49-28=21\", \"58-10=48\", \"85-2=83\", \"20+5=25\", \"30+6=36\"),\n  @(\"63+30=93\", \"55+34=89\", \"23+30=53\", \"75-17=58\", \"47-30=17\"),\n  @(\"33+28=61\", \"72+13=85\", \"31+40=71\", \"7+58=65\", \"85+2=87\"),\n  @(\"77-18=59\", \"37-9=28\", \"99-92=7\", \"84-25=59\", \"77+10=87\"),\n  @(\"64+18=82\", \"50+19=69\", \"38+2=40\", \"96-40=56\", \"79-59=20\"),\n  @(\"81-30=51\", \"81-48=33\", \"84-2=82\", \"27-16=11\", \"43+1=44\"),\n  @(\"48+34=82\", \"8+49=57\", \"35+29=64\", \"17+49=66\", \"69-55=14\"),\n  @(\"23+61=84\", \"77+14=91\", \"38+14=52\", \"95-49=46\", \"34+55=89\")\n)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n  $row = $newValues[$r]\n  for ($c = 0; $c -lt $row.Length; $c++) {\n    $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n  }\n}\n"}
